$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.827.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.665"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.69"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.97"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "37.32"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +16.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.33"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.712.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.46"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.928"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.370.37"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.787.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.63"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.37"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.26"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.66%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.91%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.24"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.134"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.17"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.81"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.63"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.67%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0280"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.65"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +18.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.89"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +10.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.43"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +10.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.108"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.07"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.202"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.25"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.66"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.15%  "
